# Adding change / Sharing a file over data inconsistancies
#
# Adds annotation notes (yellow-highlighted) next to rows that show
# discrepancies between the QCP and HumMod model columns on the
# "Hemorrhage Data" sheet, and updates the view's scroll/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Yellow fill color (BGR-encoded value Excel's Interior.Color expects for
# pure yellow, RGB 255,255,0).
$yellow = 65535

# Row 14 (Cardiac) -- O14:T14
$rng = $ws.Range("O14:T14")
$rng.Interior.Color = $yellow
$ws.Range("O14").Value = "***Cardiac Output is different between the two models"

# Row 9 (Red Cell) -- O9:U9
$rng = $ws.Range("O9:U9")
$rng.Interior.Color = $yellow
$ws.Range("O9").Value = "***Plasma volume falls again in HumMod for some unknown reason"

# Row 20 (Stroke) -- O20:T20
$rng = $ws.Range("O20:T20")
$rng.Interior.Color = $yellow
$ws.Range("O20").Value = "*** The two values differ between QCP and HumMod"

# Row 23 (Renin Act.) -- O23:T23
$rng = $ws.Range("O23:T23")
$rng.Interior.Color = $yellow
$ws.Range("O23").Value = "*** The two values differ between QCP and HumMod"

# Row 25 (Na+) -- O25:T25
$rng = $ws.Range("O25:T25")
$rng.Interior.Color = $yellow
$ws.Range("O25").Value = "*** The two values differ between QCP and HumMod"

# Row 27 (Brain Blood) -- O27:T27
$rng = $ws.Range("O27:T27")
$rng.Interior.Color = $yellow
$ws.Range("O27").Value = "*** The two values differ between QCP and HumMod"

# Row 29 (Flow) -- O29:T29
$rng = $ws.Range("O29:T29")
$rng.Interior.Color = $yellow
$ws.Range("O29").Value = "*** The two values differ between QCP and HumMod"

# Scroll the view down so row 17 is at the top, and leave the final
# selection on the last note that was added (matches the author's saved
# view state).
$ws.Activate()
$ws.Range("O29:T29").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 17
$win.ScrollColumn = 1
